# edit.ps1 - Applies the commit "Se actualizo la fecha del archivo de Word"
# 1) Wraps the English "Maximum 100 words." and "Keywords-Minimum 3 and maximum 5"
#    runs with w:proofErr spellStart/spellEnd markers (as Word's proofer would when
#    re-validating the es-CO-tagged English text), splitting the runs accordingly.
# 2) Updates the footer year from 2023 to 2024 (split into separate runs, as Word
#    does when only part of a run's text is edited in place).

$d = $word.ActiveDocument

function Replace-ParagraphWithXml($range, $partName, $partContentType, $rootTag, $innerXml) {
    $wrapper = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="' + $partName + '" pkg:contentType="' + $partContentType + '">' + `
        '<pkg:xmlData><' + $rootTag + ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
        $innerXml + '</' + $rootTag + '></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($wrapper) | Out-Null
}

# --- 1. Abstract paragraph: "Maximum 100 words." gains proofErr spell-check markers ---
$rng = $d.Content
$found = $rng.Find.Execute("Maximum 100 words.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'Maximum 100 words.' paragraph" }
$abstractParaRange = $rng.Paragraphs(1).Range
$abstractXml = '<w:p w14:paraId="72779935" w14:textId="43AD7EA2" w:rsidR="004D72B5" w:rsidRPr="00096F16" w:rsidRDefault="009303D9" w:rsidP="00972203"><w:pPr><w:pStyle w:val="Abstract"/><w:rPr><w:i/><w:iCs/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Abstract</w:t></w:r><w:r><w:t>—</w:t></w:r><w:r w:rsidR="00D44C5D"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00D44C5D" w:rsidRPr="00D44C5D"><w:t xml:space="preserve">In this space, what </w:t></w:r><w:r w:rsidR="000128F4"><w:t>was</w:t></w:r><w:r w:rsidR="00D44C5D" w:rsidRPr="00D44C5D"><w:t xml:space="preserve"> done in the development of the laboratory is described. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00D44C5D" w:rsidRPr="00096F16"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>Maximum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00D44C5D" w:rsidRPr="00096F16"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> 100 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00D44C5D" w:rsidRPr="00096F16"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>words</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00D44C5D" w:rsidRPr="00096F16"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>.</w:t></w:r></w:p>'
Replace-ParagraphWithXml $abstractParaRange "/word/document.xml" "application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" "w:document" $abstractXml

# --- 2. Keywords paragraph: "Keywords-" and "Minimum 3 and maximum 5" gain proofErr markers ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Keywords", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find 'Keywords' paragraph" }
$keywordsParaRange = $rng2.Paragraphs(1).Range
$keywordsXml = '<w:p w14:paraId="713D7512" w14:textId="48550D54" w:rsidR="009303D9" w:rsidRPr="00FC5184" w:rsidRDefault="004D72B5" w:rsidP="00972203"><w:pPr><w:pStyle w:val="Keywords"/><w:rPr><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00FC5184"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>Keywords</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00FC5184"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>—</w:t></w:r><w:r w:rsidR="00FC5184" w:rsidRPr="00FC5184"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="000128F4"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>M</w:t></w:r><w:r w:rsidR="00983011" w:rsidRPr="00983011"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>inimum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00983011" w:rsidRPr="00983011"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> 3 and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00983011" w:rsidRPr="00983011"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>maximum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00983011" w:rsidRPr="00983011"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> 5</w:t></w:r><w:r w:rsidR="00FC5184" w:rsidRPr="00FC5184"><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">.  </w:t></w:r></w:p>'
Replace-ParagraphWithXml $keywordsParaRange "/word/document.xml" "application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" "w:document" $keywordsXml

# --- 3. Footer: year "2023" -> "2024" (UFPS - 20|23 -> UFPS - 20|2|4) ---
$sec = $d.Sections(1)
$ftr = $sec.Footers(2)
$footerRange = $ftr.Range
$footerXml = '<w:p w14:paraId="7F6BE745" w14:textId="786F399C" w:rsidR="001A3B3D" w:rsidRPr="00147846" w:rsidRDefault="00D74745" w:rsidP="0056610F"><w:pPr><w:pStyle w:val="Piedepgina"/><w:jc w:val="left"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">PROGRAMA DE INGENIERÍA ELECTRÓNICA </w:t></w:r><w:r w:rsidR="00147846"><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-CO"/></w:rPr><w:t>– UFPS - 20</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-CO"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-CO"/></w:rPr><w:t>4</w:t></w:r></w:p>'
Replace-ParagraphWithXml $footerRange "/word/footer1.xml" "application/vnd.openxmlformats-officedocument.wordprocessingml.footer+xml" "w:ftr" $footerXml

Write-Output "Abstract paragraph: $($abstractParaRange.Text)"
Write-Output "Keywords paragraph: $($keywordsParaRange.Text)"
Write-Output "Footer: $($ftr.Range.Text)"
